# Revert "Create pkgdown site"
#
# Two content changes:
#   1) the cached date shown by the datetimeFigureOut footer field on the
#      slide master and on every slide layout goes from 9/30/24 back to
#      9/20/24.
#   2) the "GISCO data" label on slide 1 becomes "GSICO data".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Fix the cached date text ("9/30/24" -> "9/20/24") everywhere it shows
#    up: the slide master's Date Placeholder and the Date Placeholder on
#    each of the slide layouts.
# ---------------------------------------------------------------------
function Set-CachedDate($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "9/30/24") {
                $tr.Text = "9/20/24"
            }
        }
    }
}

$master = $p.SlideMaster
Set-CachedDate $master

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    Set-CachedDate $master.CustomLayouts.Item($li)
}

# ---------------------------------------------------------------------
# 2) Fix the "GISCO data" label on slide 1 (Rounded Rectangle 4).
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$label = $slide1.Shapes.Item("Rounded Rectangle 4")
$label.TextFrame.TextRange.Text = "GSICO data"
